$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'62.024.43"
$ws.Range("E2").Value = "'  -0.32%  "

# Row 3
$ws.Range("D3").Value = "'2.422.72"
$ws.Range("E3").Value = "'  -0.06%  "

# Row 4
$ws.Range("E4").Value = "'  +0.05%  "

# Row 5
$ws.Range("D5").Value = "'563.06"
$ws.Range("E5").Value = "'  -0.29%  "

# Row 6
$ws.Range("D6").Value = "'143.42"
$ws.Range("E6").Value = "'  -0.90%  "

# Row 7
$ws.Range("E7").Value = "'  -0.04%  "

# Row 8
$ws.Range("E8").Value = "'  -0.60%  "

# Row 9
$ws.Range("E9").Value = "'  -7.70%  "

# Row 10
$ws.Range("E10").Value = "'  -0.56%  "

# Row 11
$ws.Range("E11").Value = "'  -0.45%  "

# Row 12
$ws.Range("E12").Value = "'  -4.01%  "

# Row 13
$ws.Range("E13").Value = "'  -1.44%  "

# Row 14
$ws.Range("D14").Value = "'26.10"
$ws.Range("E14").Value = "'  +0.02%  "

# Row 15
$ws.Range("E15").Value = "'  -2.30%  "

# Row 16
$ws.Range("D16").Value = "'2.858.98"
$ws.Range("E16").Value = "'  -0.06%  "

# Row 17
$ws.Range("D17").Value = "'61.985.19"
$ws.Range("E17").Value = "'  -0.01%  "

# Row 18
$ws.Range("D18").Value = "'2.429.04"
$ws.Range("E18").Value = "'  +0.26%  "

# Row 19
$ws.Range("D19").Value = "'11.28"
$ws.Range("E19").Value = "'  -0.24%  "

# Row 20
$ws.Range("D20").Value = "'323.65"
$ws.Range("E20").Value = "'  -0.56%  "

# Row 21
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").Value = "'4.14"
$ws.Range("E21").Value = "'  -1.50%  "

# Row 22
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'6.82"
$ws.Range("E22").Value = "'  +0.66%  "

# Row 23
$ws.Range("E23").Value = "'  -0.05%  "

# Row 24
$ws.Range("D24").Value = "'67.14"
$ws.Range("E24").Value = "'  +2.44%  "

# Row 25
$ws.Range("D25").Value = "'1.73"
$ws.Range("E25").Value = "'  +0.19%  "

# Row 26
$ws.Range("E26").Value = "'  -3.08%  "

# Row 27
$ws.Range("D27").Value = "'554.71"
$ws.Range("E27").Value = "'  -5.83%  "

# Row 28
$ws.Range("D28").Value = "'2.543.36"

# Row 29
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "'  -0.03%  "

# Row 30
$ws.Range("D30").Value = "'0.0₃0930"
$ws.Range("E30").Value = "'  -2.26%  "

# Row 31
$ws.Range("D31").Value = "'8.18"
$ws.Range("E31").Value = "'  -0.84%  "

# Row 32
$ws.Range("E32").Value = "'  -4.85%  "

# Row 33
$ws.Range("E33").Value = "'  -2.18%  "

# Row 35
$ws.Range("E35").Value = "'  -3.62%  "

# Row 36
$ws.Range("E36").Value = "'  -0.01%  "

# Row 37
$ws.Range("D37").Value = "'4.74"
$ws.Range("E37").Value = "'  -1.31%  "

# Row 38
$ws.Range("E38").Value = "'  -1.20%  "

# Row 39
$ws.Range("D39").Value = "'5.47"
$ws.Range("E39").Value = "'  -4.89%  "

# Row 40
$ws.Range("D40").Value = "'153.05"
$ws.Range("E40").Value = "'  -0.73%  "

# Row 41
$ws.Range("D41").Value = "'18.62"
$ws.Range("E41").Value = "'  -0.60%  "

# Row 42
$ws.Range("D42").Value = "'1.80"
$ws.Range("E42").Value = "'  -1.94%  "

# Row 43
$ws.Range("D43").Value = "'0.996"
$ws.Range("E43").Value = "'  -0.30%  "

# Row 44
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.24"
$ws.Range("E44").Value = "'  -4.82%  "

# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'147.29"
$ws.Range("E45").Value = "'  -1.94%  "

# Row 46
$ws.Range("D46").Value = "'3.64"
$ws.Range("E46").Value = "'  -0.57%  "

# Row 47
$ws.Range("D47").Value = "'0.0529"
$ws.Range("E47").Value = "'  -2.09%  "

# Row 48
$ws.Range("E48").Value = "'  +0.42%  "

# Row 49
$ws.Range("D49").Value = "'19.84"
$ws.Range("E49").Value = "'  -3.12%  "

# Row 50
$ws.Range("D50").Value = "'0.0918"
$ws.Range("E50").Value = "'  -0.85%  "

# Row 51
$ws.Range("D51").Value = "'0.0228"
$ws.Range("E51").Value = "'  -0.68%  "
